# Updated cryptos list on Tue Jul 11 11:14:58 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) figures for
# each coin row on the active worksheet, matching the latest scrape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text such as "30.415.39" or "1.000". Several of these
# look like plain numbers to Excel's auto-detection (e.g. "1.000", "246.51"),
# which would otherwise be silently coerced into numeric values and lose
# their original formatting (trailing zeros, etc). Force the whole price
# column to Text first so every assignment below is stored verbatim.
$ws.Range("D2:D51").NumberFormat = "@"

$rows = @(
    @{ Row = 2;  Price = "30.415.39";    Volume = "  +0.74%  " },
    @{ Row = 3;  Price = "1.869.13";     Volume = "  +0.37%  " },
    @{ Row = 4;  Price = "1.000";        Volume = "  +0.05%  " },
    @{ Row = 5;  Price = "246.51";       Volume = "  +1.35%  " },
    @{ Row = 6;  Price = "1.000";        Volume = "  +0.04%  " },
    @{ Row = 7;  Price = "0.4742";       Volume = "  +0.56%  " },
    @{ Row = 8;  Price = "0.2906";       Volume = "  +1.80%  " },
    @{ Row = 9;  Price = "0.06496";      Volume = "  +0.43%  " },
    @{ Row = 10; Price = "21.99";        Volume = "  +6.07%  " },
    @{ Row = 11; Price = "0.07717";      Volume = "  +0.34%  " },
    @{ Row = 12; Price = "97.48";        Volume = "  +3.41%  " },
    @{ Row = 13; Price = "0.7382";       Volume = "  +8.19%  " },
    @{ Row = 14; Price = "1.870.71";     Volume = "  +0.40%  " },
    @{ Row = 15; Price = "5.113";        Volume = "  +0.95%  " },
    @{ Row = 16; Price = "274.04";       Volume = "  +1.90%  " },
    @{ Row = 17; Price = "30.412.02";    Volume = "  +0.74%  " },
    @{ Row = 18; Price = "13.35";        Volume = "  +0.15%  " },
    @{ Row = 19; Price = "0.000007551";  Volume = "  +0.26%  " },
    @{ Row = 20; Price = $null;          Volume = "  +0.06%  " },
    @{ Row = 21; Price = "2.117.25";     Volume = "  +0.35%  " },
    @{ Row = 22; Price = "1.000";        Volume = "  +0.06%  " },
    @{ Row = 23; Price = "5.221";        Volume = "  +0.63%  " },
    @{ Row = 24; Price = "6.163";        Volume = "  +0.96%  " },
    @{ Row = 25; Price = "9.277";        Volume = "  -0.49%  " },
    @{ Row = 26; Price = "164.21";       Volume = "  -0.83%  " },
    @{ Row = 27; Price = "18.83";        Volume = "  +0.44%  " },
    @{ Row = 28; Price = "1.929";        Volume = "  +2.20%  " },
    @{ Row = 29; Price = "0.09997";      Volume = "  +1.74%  " },
    @{ Row = 30; Price = "1.365";        Volume = "  -0.47%  " },
    @{ Row = 31; Price = "1.506";        Volume = "  -0.72%  " },
    @{ Row = 32; Price = "4.304";        Volume = "  +1.61%  " },
    @{ Row = 33; Price = "4.150";        Volume = "  +4.42%  " },
    @{ Row = 34; Price = "0.04834";      Volume = "  +2.73%  " },
    @{ Row = 35; Price = "1.119";        Volume = "  +0.90%  " },
    @{ Row = 36; Price = "0.6969";       Volume = "  +1.71%  " },
    @{ Row = 37; Price = $null;          Volume = "  +0.11%  " },
    @{ Row = 39; Price = "0.01858";      Volume = "  +0.44%  " },
    @{ Row = 40; Price = "2.745";        Volume = "  +0.68%  " },
    @{ Row = 41; Price = "6.299";        Volume = "  -1.40%  " },
    @{ Row = 42; Price = "72.84";        Volume = "  +3.41%  " },
    @{ Row = 43; Price = "1.969";        Volume = "  +4.61%  " },
    @{ Row = 44; Price = "0.4189";       Volume = "  +3.08%  " },
    @{ Row = 45; Price = $null;          Volume = "  +0.06%  " },
    @{ Row = 46; Price = "0.8361";       Volume = "  +0.07%  " },
    @{ Row = 47; Price = "102.09";       Volume = "  +0.13%  " },
    @{ Row = 48; Price = "9.231";        Volume = "  +0.91%  " },
    @{ Row = 49; Price = "7.011";        Volume = "  +1.18%  " },
    @{ Row = 51; Price = "925.42";       Volume = "  +0.32%  " }
)

foreach ($item in $rows) {
    $r = $item.Row
    if ($null -ne $item.Price) {
        $ws.Cells.Item($r, 4).Value = $item.Price
    }
    $ws.Cells.Item($r, 5).Value = $item.Volume
}
